$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# --- Cells that change data type (numeric <-> text "0"/"***.*") ---
# Use Copy() from a same-styled, unmodified reference cell to get the exact
# style index + shared-string text, then overwrite the value if needed.
$ws.Range("D14").Copy($ws.Range("C15"))   # -> t="s" s="14" "0"
$ws.Range("D14").Copy($ws.Range("D18"))   # -> t="s" s="14" "0"
$ws.Range("E14").Copy($ws.Range("E18"))   # -> t="s" s="14" "***.*"
$ws.Range("G15").Copy($ws.Range("F22"))   # -> numeric s="15"
$ws.Range("F22").Value = 1

# --- Plain numeric value updates ---
# Row 14
$ws.Range("N14").Value = -95
# Row 15
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 300
$ws.Range("L15").Value = 75
# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = -26.086956521739
$ws.Range("I16").Value = 190
$ws.Range("J16").Value = 202
$ws.Range("K16").Value = -5.940594059405
$ws.Range("L16").Value = 9.826589595375
$ws.Range("M16").Value = -19.831223628692
$ws.Range("N16").Value = -75.641025641025
# Row 17
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 33.333333333333
$ws.Range("I17").Value = 323
$ws.Range("J17").Value = 314
$ws.Range("K17").Value = 2.866242038216
$ws.Range("L17").Value = 3.525641025641
$ws.Range("M17").Value = 144.69696969697
$ws.Range("N17").Value = -14.095744680851
# Row 18
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -16.666666666666
$ws.Range("I18").Value = 124
$ws.Range("K18").Value = -10.144927536231
$ws.Range("L18").Value = 29.166666666666
$ws.Range("M18").Value = -53.383458646616
$ws.Range("N18").Value = -88.268684957426
# Row 19
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -9.090909090909
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -25.423728813559
$ws.Range("I19").Value = 553
$ws.Range("J19").Value = 597
$ws.Range("K19").Value = -7.370184254606
$ws.Range("L19").Value = 22.345132743362
$ws.Range("M19").Value = 74.447949526813
$ws.Range("N19").Value = 4.734848484848
# Row 20
$ws.Range("C20").Value = 13
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 160
$ws.Range("F20").Value = 23
$ws.Range("H20").Value = 4.545454545454
$ws.Range("I20").Value = 234
$ws.Range("J20").Value = 199
$ws.Range("K20").Value = 17.587939698492
$ws.Range("L20").Value = 33.714285714285
$ws.Range("M20").Value = -14.285714285714
$ws.Range("N20").Value = -91.978059650325
# Row 21
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 50
$ws.Range("F21").Value = 121
$ws.Range("G21").Value = 132
$ws.Range("H21").Value = -8.333333333333
$ws.Range("I21").Value = 1446
$ws.Range("J21").Value = 1473
$ws.Range("K21").Value = -1.832993890020
$ws.Range("L21").Value = 18.040816326530
$ws.Range("M21").Value = 15.958299919807
$ws.Range("N21").Value = -74.662694936043
# Row 22
$ws.Range("D22").Value = 1
$ws.Range("H22").Value = -66.666666666666
$ws.Range("I22").Value = 18
$ws.Range("J22").Value = 24
$ws.Range("K22").Value = -25
$ws.Range("L22").Value = 28.571428571428
$ws.Range("M22").Value = -28
# Row 24
$ws.Range("C24").Value = 35
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 94.444444444444
$ws.Range("F24").Value = 125
$ws.Range("G24").Value = 90
$ws.Range("H24").Value = 38.888888888888
$ws.Range("I24").Value = 1183
$ws.Range("J24").Value = 1217
$ws.Range("K24").Value = -2.793755135579
$ws.Range("L24").Value = 40
$ws.Range("M24").Value = 92.357723577235
# Row 25
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -41.666666666666
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = -19.607843137254
$ws.Range("I25").Value = 475
$ws.Range("J25").Value = 476
$ws.Range("K25").Value = -0.210084033613
$ws.Range("L25").Value = 13.636363636363
$ws.Range("M25").Value = 6.026785714285
# Row 26
$ws.Range("C26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 66.666666666666
$ws.Range("I26").Value = 32
$ws.Range("J26").Value = 35
$ws.Range("K26").Value = -8.571428571428
$ws.Range("L26").Value = 77.777777777777
# Row 27
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -25
$ws.Range("I27").Value = 42
$ws.Range("K27").Value = -20.754716981132
$ws.Range("L27").Value = 16.666666666666
